# Update the two "Student on a budget" / "Busy Parent" user stories to the
# new "image upload -> nutritional analysis / recipe recommendation" stories.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Stories")

# Row 2: Student on a budget -> track my food expense / budget  ==>
#        Student -> upload an image of food / receive video recipes (Recipe Recommendation)
# Row 3: Busy Parent -> easily upload my grocery receipts / auto-update pantry ==>
#        Parent -> see detailed nutritional info / make informed dietary choices (Nutritional Analysis)
$ws.Range("C2").Value = "Student "
$ws.Range("D2").Value = "upload an image of food"
$ws.Range("C3").Value = "Parent"
$ws.Range("D3").Value = "see detailed nutritional information for the food in the uploaded image"
$ws.Range("E3").Value = "I can make informed dietary choices"
$ws.Range("E2").Value = "I can receive video recipes related to the image"
$ws.Range("F3").Value = "Nutritional Analysis"
$ws.Range("F2").Value = "Recipe Recommendation"

# Update the sheet view: zoom level and active selection moved.
$ws.Application.ActiveWindow.Zoom = 125
$ws.Range("D21").Select()
